# TASK_38 — invite-game backlog items + promotion of "9.4" task to the
# main task sheet. Mirrors the author's commit: a batch of new backlog
# entries on "Бэклог задач" plus closing out task #38 ("9.3 ...") and
# opening task #39 ("9.4 ...") on "Задачи".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Задачи
$ws2 = $wb.Worksheets.Item(2)   # Бэклог задач

# --- 1. New backlog rows (sheet "Бэклог задач"), rows 26..33 ---------------

$backlog = @(
    @{ Row = 26; Text = "Добавить поясняющие сообщение для DB, в случае если невозможно подключиться, сообщить о конфиге и т.д."; Date = 42018.599305555559 },
    @{ Row = 27; Text = "Добавить поясняющие сообщение для WebSocketServer, в случае если невозможно подключиться, сообщить о конфиге и т.д., занятости 80 порта(skype,apache?)"; Date = 42018.599305555559 },
    @{ Row = 28; Text = "Вынести конфиг клиента в отдельный файл."; Date = 42018.599305555559 },
    @{ Row = 29; Text = "Поменять в photoInfo .id на userId. Ибо id  не фото, а юзера."; Date = 42018.623611111114 },
    @{ Row = 30; Text = "При создании игры по приглашению, если юзер оффлайн, нужно сообщить клиенту, что триндец, иначе так и будет висеть в ожидании игры."; Date = 42018.640972222223 },
    @{ Row = 31; Text = "Фотографию оппонента вывести"; Date = 42018.645138888889 },
    @{ Row = 32; Text = "Для игры по приглашению учесть знак обоих игроков, а не только одного."; Date = 42018.645833333336 },
    @{ Row = 33; Text = "Рефакторинг: разделить random, robot and invitation game."; Date = 42018.697916666664 }
)

foreach ($item in $backlog) {
    $r = $item.Row
    $ws2.Cells.Item($r, 2).Value = $item.Text
    $ws2.Cells.Item($r, 3).Value = $item.Date

    # copy the date-time number format from the row above (keeps the same
    # cell style index instead of minting a fresh one)
    $ws2.Cells.Item($r - 1, 3).Copy()
    $ws2.Cells.Item($r, 3).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# column B grew a lot wider once the long WebSocketServer sentence landed in it
$ws2.Columns.Item(2).AutoFit()

# --- 2. "Задачи" sheet: finish task 38, start task 39 -----------------------

# Task #38 (9.3 - invite accept) finished on 2015-01-14
$ws1.Cells.Item(41, 5).Value = 42018.723611111112

# Task #39 (9.4 - user state) created
$ws1.Cells.Item(42, 2).Value = "9.4 – Функционал состояния пользователя. (APIUserState)"
$ws1.Cells.Item(42, 3).Value = 42018.723611111112
$ws1.Cells.Item(42, 4).Value = 42018.723611111112

# --- 3. Selections left by the author while doing this edit ----------------

$ws2.Range("A34").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("D42").Select() | Out-Null
